$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 1151
$ws1.Range("F12").Value = 319
$ws1.Range("F16").Value = 124
$ws1.Range("F17").Value = 1272
$ws1.Range("F19").Value = 252
$ws1.Range("F26").Value = 1070
$ws1.Range("F28").Value = 3305
$ws1.Range("F31").Value = 1466

# Sheet "本地生活" (sheet 3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 775

# Sheet "全部类型" (sheet 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 775
$ws4.Range("F10").Value = 1151
$ws4.Range("F24").Value = 319
$ws4.Range("F28").Value = 124
$ws4.Range("F29").Value = 1272
$ws4.Range("F31").Value = 252
$ws4.Range("F40").Value = 1070
$ws4.Range("F42").Value = 3305
$ws4.Range("F45").Value = 1466
